$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), matching the style of the existing
# header cells (e.g. G1: bold, bordered, centered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New "Save" column values for each data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
